# Refresh crypto price / 1h-volume figures (and fix the NEARProtocol /
# PancakeSwap row ordering) per the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.823.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.91%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.443.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.14%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = '''583.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.56%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''173.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.67%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.01%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.34%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''3.441.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.12%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.131'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.03%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''6.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.79%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -2.69%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.038.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.15%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  +2.12%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''28.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -6.95%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''65.827.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.89%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -0.66%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.437.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.03%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''5.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.34%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''13.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.21%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''368.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -1.49%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''7.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.05%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''72.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +1.87%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.22%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +1.06%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +3.46%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''9.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.17%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.178'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +3.52%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.05%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''23.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -0.51%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = '''NEARProtocol'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = '''5.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -2.25%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = '''PancakeSwap'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = '''1.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -1.18%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.02%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -4.30%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -2.18%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +0.69%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''160.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +0.40%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.880'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.68%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''28.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +5.40%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.36%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -0.73%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''2.754.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.85%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''6.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -0.84%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''4.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +0.59%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.0679'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -1.77%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''40.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.00%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -2.43%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.0289'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.27%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''325.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +2.19%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -0.21%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''6.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.89%  '
$ws.Range("E51").Style = "Normal"
